$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2722  # was 2719
$ws.Range("F10").Value = 6080  # was 6074
$ws.Range("F12").Value = 68  # was 67
$ws.Range("F16").Value = 96  # was 95
$ws.Range("F17").Value = 13  # was 11
$ws.Range("F18").Value = 2552  # was 2551
$ws.Range("F22").Value = 282  # was 281
$ws.Range("F25").Value = 1017  # was 1015
$ws.Range("F28").Value = 8  # was 6
$ws.Range("F30").Value = 11  # was 9
$ws.Range("F32").Value = 287  # was 286
$ws.Range("F33").Value = 566  # was 565
$ws.Range("F38").Value = 1012  # was 1011
$ws.Range("F41").Value = 267  # was 265
$ws.Range("F43").Value = 2516  # was 2515
$ws.Range("F44").Value = 55  # was 54
$ws.Range("F48").Value = 72  # was 71

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 79  # was 78

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F9").Value = 1804  # was 1805
$ws.Range("F10").Value = 2409  # was 2408
$ws.Range("F11").Value = 803  # was 801
$ws.Range("F12").Value = 701  # was 700

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 2722  # was 2719
$ws.Range("F10").Value = 2409  # was 2408
$ws.Range("F11").Value = 6080  # was 6075
$ws.Range("F12").Value = 803  # was 801
$ws.Range("F15").Value = 68  # was 67
$ws.Range("F17").Value = 96  # was 95
$ws.Range("F18").Value = 2552  # was 2551
$ws.Range("F22").Value = 282  # was 281
$ws.Range("F30").Value = 287  # was 286
$ws.Range("F31").Value = 566  # was 565
$ws.Range("F40").Value = 267  # was 265
$ws.Range("F43").Value = 2516  # was 2515
$ws.Range("F47").Value = 72  # was 71
